$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.653.53'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '3.156.84'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '529.64'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '139.57'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.544'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +15.91%  '
$ws.Range('E9').Value = '  -0.20%  '
$ws.Range('E10').Value = '  +5.35%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.112'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.53%  '
$ws.Range('E12').Value = '  +3.39%  '
$ws.Range('D13').Value = '3.699.86'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '25.94'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('E15').Value = '  +4.54%  '
$ws.Range('D16').Value = '58.716.17'
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('E17').Value = '  +3.04%  '
$ws.Range('D18').Value = '3.155.47'
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.99'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.60%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '8.13'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '376.25'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.92%  '
$ws.Range('B22').Value = 'LEO'
$ws.Range('C22').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.80'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.13%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.532'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +4.83%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '69.59'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.58%  '
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.27'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +12.29%  '
$ws.Range('D29').Value = '0.0₃0867'
$ws.Range('E29').Value = '  -1.38%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '22.37'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +4.00%  '
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.08'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.05%  '
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '6.27'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.82%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '158.40'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('E37').Value = '  +4.21%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '24.93'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.82%  '
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').Value = '2.645.74'
$ws.Range('E40').Value = '  +5.23%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0686'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.82%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.26'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.93%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.721'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.60%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '39.10'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.40%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0288'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +7.09%  '
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').Value = '3.196.72'
$ws.Range('E47').Value = '  +0.45%  '
$ws.Range('E48').Value = '  +14.06%  '
$ws.Range('E49').Value = '  +1.68%  '
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '20.01'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.24%  '
